$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3471.4285
$ws.Range("I2").Value = 5400
$ws.Range("J2").Value = 2700
$ws.Range("K2").Value = 5400
$ws.Range("L2").Value = 2700
$ws.Range("M2").Value = -5287
$ws.Range("N2").Value = -2926
$ws.Range("H98").Value = 2830.2917
$ws.Range("I98").Value = 2402.4358
$ws.Range("K98").Value = 2402.4358
$ws.Range("M98").Value = -904.4358000000002
$ws.Range("H122").Value = 2830.2917
$ws.Range("I122").Value = 2402.4358
$ws.Range("K122").Value = 7207.307400000001
$ws.Range("M122").Value = -4757.307400000001
$ws.Range("H132").Value = 2317743.8
$ws.Range("I132").Value = 2565841.5
$ws.Range("K132").Value = 7697524.5
$ws.Range("M132").Value = -7694994.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 45000
$ws.Range("J7").Value = 45000
$ws.Range("L7").Value = 45000
$ws.Range("N7").Value = -45228
$ws.Range("H52").Value = 52499.5
$ws.Range("J52").Value = 75000
$ws.Range("L52").Value = 75000
$ws.Range("N52").Value = -75636
$ws.Range("H74").Value = 128054.06
$ws.Range("I74").Value = 137309
$ws.Range("K74").Value = 137309
$ws.Range("M74").Value = -136435
$ws.Range("H77").Value = 128054.06
$ws.Range("I77").Value = 137309
$ws.Range("K77").Value = 686545
$ws.Range("M77").Value = -682177
$ws.Range("H115").Value = 47599
$ws.Range("J115").Value = 47599
$ws.Range("L115").Value = 47599
$ws.Range("N115").Value = -50733
$ws.Range("H132").Value = 2073.2778
$ws.Range("I132").Value = 1867.225
$ws.Range("J132").Value = 2662
$ws.Range("K132").Value = 5601.674999999999
$ws.Range("L132").Value = 7986
$ws.Range("M132").Value = -3071.674999999999
$ws.Range("N132").Value = -13046

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 20826
$ws.Range("J2").Value = 20826
$ws.Range("L2").Value = 20826
$ws.Range("N2").Value = -21052
$ws.Range("H13").Value = 47300
$ws.Range("J13").Value = 47300
$ws.Range("L13").Value = 47300
$ws.Range("N13").Value = -47636
$ws.Range("H50").Value = 125000
$ws.Range("J50").Value = 125000
$ws.Range("L50").Value = 125000
$ws.Range("N50").Value = -126148
$ws.Range("H52").Value = 49917.4
$ws.Range("J52").Value = 49917.4
$ws.Range("L52").Value = 49917.4
$ws.Range("N52").Value = -50443.4
$ws.Range("H53").Value = 50000
$ws.Range("J53").Value = 50000
$ws.Range("L53").Value = 50000
$ws.Range("N53").Value = -51148
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H94").Value = 1717.5333
$ws.Range("I94").Value = 247.4
$ws.Range("J94").Value = 2452.6
$ws.Range("K94").Value = 247.4
$ws.Range("L94").Value = 2452.6
$ws.Range("M94").Value = 203.6
$ws.Range("N94").Value = -3354.6
$ws.Range("H109").Value = 44998
$ws.Range("J109").Value = 44998
$ws.Range("L109").Value = 44998
$ws.Range("N109").Value = -47772
$ws.Range("H118").Value = 49000
$ws.Range("J118").Value = 49000
$ws.Range("L118").Value = 49000
$ws.Range("N118").Value = -52314
$ws.Range("H121").Value = 49917.4
$ws.Range("J121").Value = 49917.4
$ws.Range("L121").Value = 49917.4
$ws.Range("N121").Value = -53411.4
$ws.Range("H134").Value = 2317.2341
$ws.Range("I134").Value = 2093.4285
$ws.Range("J134").Value = 4197.2
$ws.Range("K134").Value = 6280.2855
$ws.Range("L134").Value = 12591.6
$ws.Range("M134").Value = -3745.2855
$ws.Range("N134").Value = -17661.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5001610.5
$ws.Range("I31").Value = 5264590
$ws.Range("K31").Value = 5264590
$ws.Range("M31").Value = -5264295
$ws.Range("H34").Value = 5001610.5
$ws.Range("I34").Value = 5264590
$ws.Range("K34").Value = 5264590
$ws.Range("M34").Value = -5264388
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H109").Value = 39998
$ws.Range("J109").Value = 39998
$ws.Range("L109").Value = 39998
$ws.Range("N109").Value = -42078
$ws.Range("H114").Value = 24841
$ws.Range("J114").Value = 4684
$ws.Range("L114").Value = 4684
$ws.Range("N114").Value = -13362
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H120").Value = 21991.6
$ws.Range("J120").Value = 21991.6
$ws.Range("L120").Value = 21991.6
$ws.Range("N120").Value = -29249.6
$ws.Range("H121").Value = 38499
$ws.Range("J121").Value = 38499
$ws.Range("L121").Value = 38499
$ws.Range("N121").Value = -41119
$ws.Range("H132").Value = 20256.031
$ws.Range("I132").Value = 25088.16
$ws.Range("J132").Value = 2998.4285
$ws.Range("K132").Value = 75264.48
$ws.Range("L132").Value = 8995.2855
$ws.Range("M132").Value = -72734.48
$ws.Range("N132").Value = -14055.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2124.8572
$ws.Range("I113").Value = 720.25
$ws.Range("J113").Value = 3997.6667
$ws.Range("K113").Value = 2160.75
$ws.Range("L113").Value = 11993.0001
$ws.Range("M113").Value = 9.25
$ws.Range("N113").Value = -16333.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H45").Value = 19378.428
$ws.Range("J45").Value = 25883
$ws.Range("L45").Value = 25883
$ws.Range("N45").Value = -27001
$ws.Range("H51").Value = 47500
$ws.Range("J51").Value = 47500
$ws.Range("L51").Value = 47500
$ws.Range("N51").Value = -48518
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H126").Value = 5253.5713
$ws.Range("I126").Value = 2616.6667
$ws.Range("K126").Value = 7850.000100000001
$ws.Range("M126").Value = -5380.000100000001
$ws.Range("H132").Value = 1868.7906
$ws.Range("I132").Value = 1819.9474
$ws.Range("J132").Value = 2240
$ws.Range("K132").Value = 5459.8422
$ws.Range("L132").Value = 6720
$ws.Range("M132").Value = -2929.8422
$ws.Range("N132").Value = -11780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 34996.668
$ws.Range("J109").Value = 34996.668
$ws.Range("L109").Value = 34996.668
$ws.Range("N109").Value = -37770.668
$ws.Range("H132").Value = 2900.6316
$ws.Range("I132").Value = 3053.647
$ws.Range("J132").Value = 1600
$ws.Range("K132").Value = 9160.940999999999
$ws.Range("L132").Value = 4800
$ws.Range("M132").Value = -6630.940999999999
$ws.Range("N132").Value = -9860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2487
$ws.Range("I96").Value = 900
$ws.Range("J96").Value = 3016
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 3016
$ws.Range("M96").Value = 473
$ws.Range("N96").Value = -5762
$ws.Range("H109").Value = 79985
$ws.Range("J109").Value = 79985
$ws.Range("L109").Value = 79985
$ws.Range("N109").Value = -82759
$ws.Range("H122").Value = 11394484
$ws.Range("I122").Value = 11659007
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 34977021
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -34974571
$ws.Range("H136").Value = 15335.02
$ws.Range("J136").Value = 4199.75
$ws.Range("L136").Value = 12599.25
$ws.Range("N136").Value = -17699.25
